$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price, Volume) to Text format so numeric-looking
# strings (e.g. "30.665.15", "1.000", "248.41") are preserved verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.665.15"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.959.25"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "248.41"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.4813"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("D8").Value = "44.62"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "0.2925"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "0.06755"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "108.85"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "19.13"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.961.20"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").Value = "0.07743"
$ws.Range("D15").Value = "5.469"
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("D16").Value = "0.6932"
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("D17").Value = "291.83"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "30.678.57"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "5.669"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").Value = "13.14"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "0.000007692"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.217.58"
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "6.602"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.907"
$ws.Range("E26").Value = "  +4.82%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "169.53"
$ws.Range("E27").Value = "  +2.94%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "19.98"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.171"
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.1066"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "1.443"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.857"
$ws.Range("E32").Value = "  +19.96%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.467"
$ws.Range("E33").Value = "  +7.92%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.05093"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7697"
$ws.Range("E35").Value = "  +4.05%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.177"
$ws.Range("E36").Value = "  +3.74%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02039"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").Value = "2.721"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.716"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.489"
$ws.Range("E40").Value = "  +11.21%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "2.124"
$ws.Range("E41").Value = "  +5.59%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "109.99"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8844"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4460"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "69.91"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9993"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.459"
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1276"
$ws.Range("E48").Value = "  +3.62%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.399"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "36.04"
$ws.Range("E50").Value = "  +3.31%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "933.05"
$ws.Range("E51").Value = "  +9.99%  "
